$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "criterion_or_website"
$ws.Range("B1").Value = "website_1_value"
$ws.Range("C1").Value = "website_2_value"
$ws.Range("D1").Value = "website_3_value"
$ws.Range("E1").Value = "website_4_value"
$ws.Range("F1").Value = "conclusion"
$ws.Range("B2").Value = "Utiliza fuentes sans-serif claras y legibles, como Segoe UI, con una jerarquía consistente para títulos y texto principal."
$ws.Range("C2").Value = "Emplea tipografías sans-serif modernas y limpias, priorizando la legibilidad en subtítulos y descripciones de fotos."
$ws.Range("D2").Value = "Fuentes sans-serif profesionales como Open Sans, garantizando alta legibilidad para contenido empresarial y perfiles detallados."
$ws.Range("E2").Value = "Tipografía sans-serif estándar y sencilla, lo que asegura buena legibilidad para nombres de productos y títulos."
$ws.Range("F2").Value = "Podría explorar una tipografía más distintiva o una paleta de fuentes complementarias para añadir carácter y diferenciación visual."
$ws.Range("B3").Value = "Predominan azules y blancos, transmitiendo confianza y familiaridad, con un contraste elevado para la lectura."
$ws.Range("C3").Value = "Logo con gradientes vibrantes, interfaz principalmente blanca y texto negro, destacando el contenido visual."
$ws.Range("D3").Value = "Esquema de colores corporativos: azul, blanco y gris, que proyecta profesionalismo y un entorno de trabajo serio."
$ws.Range("E3").Value = "Dominante cabecera azul brillante con texto blanco, cuerpo de la página blanco con texto negro, colores funcionales."
$ws.Range("F3").Value = "La paleta de colores es básica; considerar introducir un color secundario o acentos para destacar elementos clave y mejorar la estética general."
$ws.Range("B4").Value = "Mezcla de tono, más informal para interacciones personales y algo formal para noticias y páginas oficiales."
$ws.Range("C4").Value = "Predominantemente informal y visual, centrado en la expresión creativa y estilos de vida personales."
$ws.Range("D4").Value = "Estrictamente formal y profesional, diseñado para networking, búsqueda de empleo y desarrollo de negocios."
$ws.Range("E4").Value = "Tono formal y funcional, típico de sitios de comercio electrónico y utilidades, enfocado en información de productos."
$ws.Range("F4").Value = "El tono es adecuado para una tienda de componentes; sin embargo, un toque sutil de calidez podría mejorar la conexión con el usuario."
$ws.Range("B5").Value = "Íconos universales de redes sociales, botones de 'me gusta', perfiles de usuario y emoticonos reconocibles."
$ws.Range("C5").Value = "Íconos minimalistas (corazón, burbuja de diálogo, avión de papel), con énfasis en las imágenes de perfil de usuario."
$ws.Range("D5").Value = "Íconos profesionales para conexiones, empleos y mensajes, junto con logotipos de empresas y fotos de perfil."
$ws.Range("E5").Value = "Íconos simples y funcionales: búsqueda, perfil de usuario y un logo de ratón 'PMCity' claro y directo en la cabecera."
$ws.Range("F5").Value = "Los íconos son funcionales pero sencillos; mejorar el diseño del logo y considerar íconos más modernos y coherentes visualmente."
$ws.Range("B6").Value = "Generalmente buena accesibilidad con texto alternativo para imágenes y navegación por teclado bien implementada."
$ws.Range("C6").Value = "Accesibilidad básica, incluyendo texto alternativo para fotos y soporte para lectores de pantalla en ciertas funciones."
$ws.Range("D6").Value = "Buenas prácticas de accesibilidad, estructura semántica clara y navegación por teclado robusta para profesionales."
$ws.Range("E6").Value = "Aparentemente buen contraste de texto y navegación clara. No se observan barreras obvias en la interfaz visible."
$ws.Range("F6").Value = "Asegurar que todos los elementos interactivos sean navegables por teclado y que las imágenes de productos incluyan descripciones de texto alternativo."
$ws.Range("B7").Value = "Barra superior con secciones (Inicio, Watch, Marketplace), barra lateral extensa y barra de búsqueda prominentemente."
$ws.Range("C7").Value = "Barra de navegación inferior intuitiva (Inicio, Buscar, Reels, Tienda, Perfil) y botones superiores para mensajes."
$ws.Range("D7").Value = "Navegación superior global (Inicio, Mi red, Empleos) y un menú de perfil claro, estructurado para uso profesional."
$ws.Range("E7").Value = "Cabecera con enlaces 'Información', 'Arma tu PC', 'Comparar' y una barra de búsqueda visible, navegación directa."
$ws.Range("F7").Value = "La navegación principal es clara; se podría añadir un 'carrito' de compras o un historial de búsquedas para mejorar la experiencia del usuario."
$ws.Range("B8").Value = "Organización basada en un 'feed' de noticias, secciones para grupos, páginas y perfiles personales."
$ws.Range("C8").Value = "Perfiles en cuadrícula, 'feed' y 'stories'; organización visual que enfatiza la disposición y el descubrimiento de contenido."
$ws.Range("D8").Value = "Centrado en el perfil, 'feed' de noticias, listados de empleo y páginas de empresa, con secciones estructuradas para datos profesionales."
$ws.Range("E8").Value = "Cabecera, título principal y sección de 'Componentes populares' en un diseño de cuadrícula claro y simple para productos."
$ws.Range("F8").Value = "La organización es lógica para productos; considerar filtros de búsqueda avanzada y categorización de componentes para facilitar la exploración."
$ws.Range("B9").Value = "Messenger, Marketplace, Grupos, Eventos, Juegos y transmisiones en vivo, ofreciendo un amplio conjunto de funciones."
$ws.Range("C9").Value = "Reels, Stories, IGTV, Tienda y Live, un conjunto robusto de funciones para compartir y consumir contenido multimedia."
$ws.Range("D9").Value = "Búsqueda de empleo, LinkedIn Learning, artículos profesionales y grupos temáticos, herramientas profesionales completas."
# E9's target text starts with a literal apostrophe. Assigning a value whose
# first character is "'" is treated as Excel's quote-prefix (text) marker and
# that leading character gets stripped from the stored text. Prepending an
# extra "'" means only the marker is consumed, leaving the literal leading
# apostrophe intact in the cell text; the style is reset below so no stray
# quote-prefix formatting is left behind.
$ws.Range("E9").Value = "''Arma tu PC' y 'Comparar', funciones muy relevantes para su propósito, que añaden valor directo al usuario."
$ws.Range("F9").Value = "Las funciones de 'Arma tu PC' y 'Comparar' son excelentes; se podría integrar un sistema de reseñas de usuarios o foros para productos."
$ws.Range("B10").Value = "Centro de ayuda extenso con guías detalladas para funciones, consejos de seguridad y normas comunitarias."
$ws.Range("C10").Value = "Centro de ayuda, guías para usar funciones como Reels o Stories, y configuración de privacidad."
$ws.Range("D10").Value = "LinkedIn Learning (premium), extenso centro de ayuda y artículos de consejos profesionales y guías de uso."
$ws.Range("E10").Value = "La sección 'Aprende' sugiere contenido educativo, y la función 'Arma tu PC' actúa como un tutorial guiado de construcción."
$ws.Range("F10").Value = "Sería beneficioso incluir tutoriales de montaje o guías detalladas para principiantes, accesibles directamente desde las secciones 'Aprende' o 'Información'."
$ws.Range("A11").Value = "Overall User Experience"
$ws.Range("B11").Value = "Experiencia integral, pero puede ser abrumadora por la cantidad de funciones, aunque intuitiva en lo básico."
$ws.Range("C11").Value = "Muy atractiva y visual, intuitiva para compartir medios, pero puede generar sobrecarga de contenido visual."
$ws.Range("D11").Value = "Enfocada y eficiente para uso profesional, con rutas claras para networking y búsqueda de empleo, a veces algo formal."
$ws.Range("E11").Value = "Directa y funcional para encontrar componentes de PC, con un diseño limpio y sencillo para navegar fácilmente por los productos."
$ws.Range("F11").Value = "Para mejorar la experiencia general, optimizar los tiempos de carga de imágenes y añadir animaciones sutiles a las interacciones clave."

# Restore E9's default (unstyled) appearance after the quote-prefix workaround above.
$ws.Range("E9").Style = "Normal"
